# Update with Correct Forecast output
#
# 1) "Forecast Comparison" sheet: insert a new "Week_Start_Date" column
#    (B) shifting ASIN..is_holiday_week one column to the right, change
#    the Week labels from zero-padded ("W01") to unpadded ("W1"), fill in
#    the new per-week start dates, refresh the MyForecast figures, and
#    mark is_holiday_week as a boolean.
# 2) "Summary" sheet: refresh the forecast totals/extremes that depend on
#    the corrected MyForecast numbers.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: make room for the new column ---
$ws1.Columns.Item(2).Insert()
$ws1.Range("B1").Value = "Week_Start_Date"

# Keep the week-start dates as plain text (not auto-converted to date serials)
$ws1.Range("B2:B17").NumberFormat = "@"

$weeks      = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$startDates = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")
$myForecast = @(675,649,656,651,598,485,462,452,563,508,460,474,479,432,427,441)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws1.Range("A$row").Value = $weeks[$i]
    $ws1.Range("B$row").Value = $startDates[$i]
    $ws1.Range("D$row").Value = $myForecast[$i]
    $ws1.Range("J$row").Value = $false
}

# --- Summary: refresh the dependent totals/extremes ---
$summaryRows = @(9,10,11,12,14)
foreach ($r in $summaryRows) {
    $ws2.Range("B$r").NumberFormat = "@"
}
$ws2.Range("B9").Value  = "8412"
$ws2.Range("B10").Value = "4629"
$ws2.Range("B11").Value = "2631"
$ws2.Range("B12").Value = "675"
$ws2.Range("B14").Value = "427"
